$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column D ("number of acceptance criteria") for each user story row.
$values = @{
    2  = 1
    4  = 1
    6  = 3
    8  = 1
    10 = 2
    13 = 2
    15 = 2
    17 = 3
    19 = 3
    20 = 3
    22 = 2
    24 = 2
    26 = 1
    28 = 3
    31 = 2
    33 = 1
    35 = 3
    38 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# Update the active view/selection: scroll back to the top and select C22.
$window = $excel.ActiveWindow
$window.TopLeftCell = $ws.Range("A1")
$ws.Range("C22").Select()
